$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data for an imported "egresado" (graduate) record
$ws.Range("A5").Value = "2009-89034"
$ws.Range("C5").Value = "Randall"
$ws.Range("D5").Value = "Rodriguez"
$ws.Range("E5").Value = "M"
$ws.Range("G5").Value = "randall@gmail.com"
$ws.Range("J5").Value = "adasddasdsa"
$ws.Range("K5").Value = "Ingeniería de sistemas"

# Add hyperlink for the new email address, mirroring existing rows
$ws.Hyperlinks.Add($ws.Range("G5"), "mailto:randall@gmail.com", [Type]::Missing, [Type]::Missing, "randall@gmail.com")

# Update the active cell selection
$ws.Range("C14").Select()
